$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Screws task (row 50, "Put screws on wheels") marked done, plus a couple
# of related re-scored subtasks. Dependent subtotal formulas (D30/D41/D47/D51)
# recalculate automatically.
$ws.Range("D31").Value = 80
$ws.Range("D32").Value = 100
$ws.Range("D45").Value = 90
$ws.Range("D50").Value = 100
$ws.Range("D53").Value = 80

# Update the saved selection (view state) to reflect where the author left off.
$ws.Range("N42").Select()
